# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the zh-cn and de-de report sheets to reflect the regenerated
# handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 02:59:14"
$wsZhCn.Range("H2").Value = "2016-03-22 02:59:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 02:59:18"
$wsDeDe.Range("H2").Value = "2016-03-22 02:59:41"
